$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# --- Style fixups for cells changing between text-placeholder and numeric formats ---
$ws.Range("C14").NumberFormat = '#,##0'
$ws.Range("F14").NumberFormat = '#,##0'
$ws.Range("I14").NumberFormat = '#,##0'
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("F28").NumberFormat = '#,##0'
$ws.Range("I28").NumberFormat = '#,##0'
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("F29").NumberFormat = '#,##0'
$ws.Range("I29").NumberFormat = '#,##0'
$ws.Range("C15").NumberFormat = 'General'
$ws.Range("C26").NumberFormat = 'General'

# --- Data value updates ---
# Row 14
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -85.714285714285
# Row 15
$ws.Range("C15").Value = "0"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = -50
$ws.Range("N15").Value = -75
# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 5.555555555555
$ws.Range("I16").Value = 36
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = 9.090909090909
$ws.Range("L16").Value = -20
$ws.Range("M16").Value = -36.842105263157
$ws.Range("N16").Value = -91.304347826087
# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 63
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = -13.698630136986
$ws.Range("L17").Value = -4.545454545454
$ws.Range("M17").Value = -20.253164556962
$ws.Range("N17").Value = -68.341708542713
# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = -5
$ws.Range("L18").Value = -2.564102564102
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -83.898305084745
# Row 19
$ws.Range("C19").Value = 8
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 33
$ws.Range("H19").Value = 17.857142857142
$ws.Range("I19").Value = 80
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = 14.285714285714
$ws.Range("L19").Value = 29.032258064516
$ws.Range("M19").Value = 15.942028985507
$ws.Range("N19").Value = -43.262411347517
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = 23.529411764705
$ws.Range("L20").Value = 23.529411764705
$ws.Range("M20").Value = -16
$ws.Range("N20").Value = -83.59375
# Row 21
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -13.793103448275
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = -0.9900990099
$ws.Range("I21").Value = 243
$ws.Range("J21").Value = 242
$ws.Range("K21").Value = 0.413223140495
$ws.Range("L21").Value = 3.846153846153
$ws.Range("M21").Value = -22.364217252396
$ws.Range("N21").Value = -78.702892199824
# Row 22
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = 6
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -71.428571428571
# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("I23").Value = 40
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = 53.846153846153
# Row 24
$ws.Range("C24").Value = 60
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 160.869565217391
$ws.Range("F24").Value = 153
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = 33.043478260869
$ws.Range("I24").Value = 331
$ws.Range("J24").Value = 281
$ws.Range("K24").Value = 17.793594306049
$ws.Range("L24").Value = 62.254901960784
$ws.Range("M24").Value = 109.493670886076
# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -31.25
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = -24.074074074074
$ws.Range("I25").Value = 120
$ws.Range("J25").Value = 130
$ws.Range("K25").Value = -7.692307692307
$ws.Range("L25").Value = 55.844155844155
$ws.Range("M25").Value = -14.893617021276
# Row 26
$ws.Range("C26").Value = "0"
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -80
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = -36.363636363636
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 5
$ws.Range("L27").Value = 66.666666666666
# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("I28").Value = 1
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -75
$ws.Range("M28").Value = -85.714285714285
$ws.Range("N28").Value = -97.777777777777
# Row 29
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 1
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -75
$ws.Range("M29").Value = -83.333333333333
$ws.Range("N29").Value = -97.368421052631
